# Update column G ("K") values for rows 2-16 on Sheet1.
# The repository regenerated this save_data sheet to compute the "K"
# column (strikes) differently, producing new (smaller) values.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$newValues = @{
    2  = 7
    3  = 3
    4  = 0
    5  = 0
    6  = 5
    7  = 7
    8  = 5
    9  = 3
    10 = 5
    11 = 10
    12 = 5
    13 = 2
    14 = 2
    15 = 1
    16 = 4
}

foreach ($row in $newValues.Keys) {
    $ws.Range("G$row").Value = $newValues[$row]
}
